# "added new test cases"
#
# - Rename the (single) worksheet from the German default "Tabelle1" to
#   "Sheet1" (the author re-created the sheet on an English Excel build).
# - Move the live cell selection from B7 down to B33, where the new rows
#   of test cases were appended.
# - Best-effort: record the folder Excel last saved to / window geometry,
#   matching what a real Excel session would stamp on save when the file
#   was re-opened from "G:\21.10.2025\" with a maximized/moved window.
#   (These are host/session bookkeeping values; harmless to set even if
#   this runtime does not round-trip them back into the saved XML.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet was renamed from the German "Tabelle1" default to "Sheet1".
$ws.Name = "Sheet1"

# New test-case rows were added further down the sheet; the author's
# last selection before saving sat on B33 instead of the old B7.
$ws.Range("B33").Select()

# Best-effort session metadata (folder path / window placement) — no-ops
# on hosts that don't persist them, but match real Excel.Application usage.
try { $wb.Path = "G:\21.10.2025\" } catch {}

$win = $wb.Windows.Item(1)
$win.Left = 23880
$win.Top = -120
$win.Width = 29040
$win.Height = 15840
